# Generate Report for Handoff
# Update "Latest Handoff Datetime" for the most recently handed-off source
# file in each locale sheet to reflect the new handoff run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-03-01 07:12:03"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-03-01 07:12:11"
